$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56; all rows from 56 downward shift down by one.
$ws.Rows("56:56").Insert()

# Populate the newly inserted row 56 with the new weekly record.
$ws.Cells.Item(56, 1).Value = 7
$ws.Cells.Item(56, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(56, 3).Value = "Ñuble"
$ws.Cells.Item(56, 4).Value = 44495
$ws.Cells.Item(56, 5).Value = 16
$ws.Cells.Item(56, 6).Value = 100112003
$ws.Cells.Item(56, 7).Value = "Ajo"
$ws.Cells.Item(56, 8).Value = "Chino"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 100
$ws.Cells.Item(56, 11).Value = 16000
$ws.Cells.Item(56, 12).Value = 17000
$ws.Cells.Item(56, 13).Value = 16500
$ws.Cells.Item(56, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(56, 15).Value = "China"
$ws.Cells.Item(56, 16).Value = 1650
$ws.Cells.Item(56, 17).Value = 10
$ws.Cells.Item(56, 18).Value = "Hortaliza"

# Keep the date column using the same date display format as the rest of column D.
$ws.Cells.Item(56, 4).NumberFormat = $ws.Cells.Item(57, 4).NumberFormat

Write-Output "done"
